$d = $word.ActiveDocument

# 1) The three "Pedir ..." bullet paragraphs were each split across two
#    runs (e.g. "Pedir " + "el nombre completo del paciente"). The edit
#    merges each pair back into a single run with the same combined text.
#    Doing a Find & Replace over the already-correct full text causes the
#    engine to re-emit the paragraph as one run, exactly like the diff.
$d.Content.Find.Execute("Pedir el nombre completo del paciente", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Pedir el nombre completo del paciente", 2) | Out-Null

$d.Content.Find.Execute("Pedir la fecha de nacimiento del paciente", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Pedir la fecha de nacimiento del paciente", 2) | Out-Null

$d.Content.Find.Execute("Pedir la dirección de la clínica médica ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Pedir la dirección de la clínica médica ", 2) | Out-Null

# 2) The last bullet, "Escribir datos pedidos", gets the word "pedidos"
#    wrapped in a grammar-check marker pair (<w:proofErr w:type="gramStart"/>
#    ... <w:proofErr w:type="gramEnd"/>), splitting the single run into
#    "Escribir datos " + "pedidos" with the proofErr markers around the
#    second run. Locate that exact word range and replace it in place
#    with equivalent OOXML (same text, now flanked by proofErr markers).
$last = $d.Paragraphs.Last.Range
$paraStart = $last.Start
$paraTextIdx = $last.Text.IndexOf("pedidos")
$wordStart = $paraStart + $paraTextIdx
$wordEnd = $wordStart + "pedidos".Length
$target = $d.Range($wordStart, $wordEnd)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>pedidos</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml) | Out-Null
